$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# xlShiftDown = -4121, xlFormatFromLeftOrAbove = -4163
$xlShiftDown = -4121
$xlFormatFromLeftOrAbove = -4163

# Add a new test point row 16 ("TRANS7"), cloning row 15's formatting/height
# (copy+insert keeps the existing style index instead of minting a new one)
$ws.Rows(15).Copy()
$ws.Rows(16).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A16").Value = "TRANS7"
$ws.Rows(16).RowHeight = 15

# Add another test point row 17 ("TRANS8"), cloning the row just created
$ws.Rows(16).Copy()
$ws.Rows(17).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A17").Value = "TRANS8"
$ws.Rows(17).RowHeight = 15

$ws.Range("B15").Select()
